$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 'MetaDslx.CodeAnalysis.Common\SpecialTypes.cs'
$ws.Range("B6").Value = 'internal static class SpecialTypes'
$ws.Range("C6").Value = 'public static class SpecialTypes'
$ws.Range("A7").Value = 'MetaDslx.CodeAnalysis.Common\SpecialTypes.cs'
$ws.Range("B7").Value = 'public static SpecialType GetTypeFromMetadataName'
$ws.Range("C7").Value = 'internal static SpecialType GetTypeFromMetadataName'
$ws.Range("A8").Value = 'MetaDslx.CodeAnalysis.Common\SpecialTypes.cs'
$ws.Range("B8").Value = 'public static Microsoft.Cci.PrimitiveTypeCode GetTypeCode'
$ws.Range("C8").Value = 'internal static Microsoft.Cci.PrimitiveTypeCode GetTypeCode'
$ws.Range("A9").Value = 'MetaDslx.CodeAnalysis.Common\Collections\ArrayElement.cs'
$ws.Range("B9").Value = 'internal struct ArrayElement'
$ws.Range("C9").Value = 'public struct ArrayElement'
$ws.Range("A10").Value = 'MetaDslx.CodeAnalysis.Common\Collections\ArrayElement.cs'
$ws.Range("B10").Value = 'internal T Value;'
$ws.Range("C10").Value = 'public T Value;'
$ws.Range("A11").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\Hash.cs'
$ws.Range("B11").Value = 'internal static class Hash'
$ws.Range("C11").Value = 'public static class Hash'
$ws.Range("A12").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\Hash.cs'
$ws.Range("B12").Value = 'internal static int Combine'
$ws.Range("C12").Value = 'public static int Combine'
$ws.Range("A13").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\Hash.cs'
$ws.Range("B13").Value = 'internal static int CombineValues'
$ws.Range("C13").Value = 'public static int CombineValues'
$ws.Range("A14").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\ConsList`1.cs'
$ws.Range("B14").Value = 'internal class ConsList'
$ws.Range("C14").Value = 'public class ConsList'
$ws.Range("A15").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\ConsList`1.cs'
$ws.Range("B15").Value = 'internal struct Enumerator'
$ws.Range("C15").Value = 'public struct Enumerator'
$ws.Range("A16").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\ThreeState.cs'
$ws.Range("B16").Value = 'internal enum ThreeState'
$ws.Range("C16").Value = 'public enum ThreeState'
$ws.Range("A17").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\ThreeState.cs'
$ws.Range("B17").Value = 'internal static class ThreeStateHelpers'
$ws.Range("C17").Value = 'public static class ThreeStateHelpers'
$ws.Range("A18").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\ExceptionUtilities.cs'
$ws.Range("B18").Value = 'internal'
$ws.Range("C18").Value = 'public'
$ws.Range("A19").Value = 'MetaDslx.CodeAnalysis.Common\InternalUtilities\ImmutableArrayExtensions.cs'
$ws.Range("B19").Value = 'internal'
$ws.Range("C19").Value = 'public'
$ws.Range("A20").Value = 'MetaDslx.CodeAnalysis.Common\Serialization\ObjectReader.cs'
$ws.Range("B20").Value = 'internal sealed partial class ObjectReader'
$ws.Range("C20").Value = 'public sealed partial class ObjectReader'
$ws.Range("A21").Value = 'MetaDslx.CodeAnalysis.Common\Serialization\ObjectWriter.cs'
$ws.Range("B21").Value = 'internal sealed partial class ObjectWriter'
$ws.Range("C21").Value = 'public sealed partial class ObjectWriter'
$ws.Range("A22").Value = 'MetaDslx.CodeAnalysis.Common\Serialization\IObjectWritable.cs'
$ws.Range("B22").Value = 'internal interface IObjectWritable'
$ws.Range("C22").Value = 'public interface IObjectWritable'
$ws.Range("A23").Value = 'MetaDslx.CodeAnalysis.Common\Serialization\ObjectBinder.cs'
$ws.Range("B23").Value = 'internal static class ObjectBinder'
$ws.Range("C23").Value = 'public static class ObjectBinder'
$ws.Range("A24").Value = 'MetaDslx.CodeAnalysis.Common\Serialization\ObjectBinderSnapshot.cs'
$ws.Range("B24").Value = 'internal readonly struct ObjectBinderSnapshot'
$ws.Range("C24").Value = 'public readonly struct ObjectBinderSnapshot'
$ws.Range("A25").Value = 'MetaDslx.CodeAnalysis.Common\Diagnostic\SourceLocation.cs'
$ws.Range("B25").Value = 'internal sealed class SourceLocation'
$ws.Range("C25").Value = 'public sealed class SourceLocation'
$ws.Range("A26").Value = 'MetaDslx.CodeAnalysis.Common\Diagnostic\DiagnosticBag.cs'
$ws.Range("B26").Value = 'internal class DiagnosticBag'
$ws.Range("C26").Value = 'public class DiagnosticBag'
$ws.Range("A27").Value = 'MetaDslx.CodeAnalysis.Common\Diagnostic\DiagnosticBag.cs'
$ws.Range("B27").Value = 'internal static DiagnosticBag GetInstance()'
$ws.Range("C27").Value = 'public static DiagnosticBag GetInstance()'
$ws.Range("A28").Value = 'MetaDslx.CodeAnalysis.Common\Diagnostic\DiagnosticBag.cs'
$ws.Range("B28").Value = 'internal void Free()'
$ws.Range("C28").Value = 'public void Free()'
$ws.Range("A29").Value = 'MetaDslx.CodeAnalysis.Common\Diagnostic\DiagnosticInfo.cs'
$ws.Range("B29").Value = 'internal class DiagnosticInfo'
$ws.Range("C29").Value = 'public class DiagnosticInfo'
$ws.Range("A30").Value = 'MetaDslx.CodeAnalysis.Common\Diagnostic\CommonMessageProvider.cs'
$ws.Range("B30").Value = 'internal abstract class CommonMessageProvider'
$ws.Range("C30").Value = 'public abstract class CommonMessageProvider'
$ws.Range("A31").Value = 'MetaDslx.CodeAnalysis.Common\DiagnosticAnalyzer\AsyncQueue.cs'
$ws.Range("B31").Value = 'internal sealed class AsyncQueue'
$ws.Range("C31").Value = 'public sealed class AsyncQueue'
$ws.Range("A32").Value = 'MetaDslx.CodeAnalysis.Common\DiagnosticAnalyzer\CompilationEvent.cs'
$ws.Range("B32").Value = 'internal abstract class CompilationEvent'
$ws.Range("C32").Value = 'public abstract class CompilationEvent'
$ws.Range("A33").Value = 'MetaDslx.CodeAnalysis.Common\DiagnosticAnalyzer\CompilationEvent.cs'
$ws.Range("B33").Value = 'internal CompilationEvent'
$ws.Range("C33").Value = 'public CompilationEvent'
$ws.Range("A34").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\GreenNode.cs'
$ws.Range("B34").Value = 'internal abstract class GreenNode'
$ws.Range("C34").Value = 'public abstract class GreenNode'
$ws.Range("A35").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\GreenNode.cs'
$ws.Range("B35").Value = 'internal enum NodeFlags'
$ws.Range("C35").Value = 'internal protected enum NodeFlags'
$ws.Range("A36").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\GreenNodeExtensions.cs'
$ws.Range("B36").Value = 'internal static class GreenNodeExtensions'
$ws.Range("C36").Value = 'public static class GreenNodeExtensions'
$ws.Range("A37").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxNodeOrToken.cs'
$ws.Range("B37").Value = 'internal int Position => _position;'
$ws.Range("C37").Value = 'internal int Position => _position;`n`n    public bool IsNull => _nodeOrParent is null && _token is null;'
$ws.Range("A38").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxList`1.cs'
$ws.Range("B38").Value = 'internal SyntaxList'
$ws.Range("C38").Value = 'public SyntaxList'
$ws.Range("A39").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxList`1.cs'
$ws.Range("B39").Value = 'internal SyntaxNode? Node'
$ws.Range("C39").Value = 'public SyntaxNode? Node'
$ws.Range("A40").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SeparatedSyntaxList.cs'
$ws.Range("B40").Value = 'internal SeparatedSyntaxList'
$ws.Range("C40").Value = 'public SeparatedSyntaxList'
$ws.Range("A41").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SeparatedSyntaxList.cs'
$ws.Range("B41").Value = 'internal SyntaxNode? Node'
$ws.Range("C41").Value = 'public SyntaxNode? Node'
$ws.Range("A42").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxToken.cs'
$ws.Range("B42").Value = 'internal SyntaxToken'
$ws.Range("C42").Value = 'public SyntaxToken'
$ws.Range("A43").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxToken.cs'
$ws.Range("B43").Value = 'internal GreenNode? Node'
$ws.Range("C43").Value = 'public GreenNode? Node'
$ws.Range("A44").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxTokenList.cs'
$ws.Range("B44").Value = 'internal SyntaxTokenList'
$ws.Range("C44").Value = 'public SyntaxTokenList'
$ws.Range("A45").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxTokenList.cs'
$ws.Range("B45").Value = 'internal GreenNode? Node'
$ws.Range("C45").Value = 'public GreenNode? Node'
$ws.Range("A46").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxTrivia.cs'
$ws.Range("B46").Value = 'internal SyntaxTrivia'
$ws.Range("C46").Value = 'public SyntaxTrivia'
$ws.Range("A47").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxTree.cs'
$ws.Range("B47").Value = 'internal virtual bool SupportsLocations'
$ws.Range("C47").Value = 'public virtual bool SupportsLocations'
$ws.Range("A48").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\SyntaxReference.cs'
$ws.Range("B48").Value = 'internal Location GetLocation'
$ws.Range("C48").Value = 'public Location GetLocation'
$ws.Range("A49").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\GreenNodeExtensions.cs'
$ws.Range("B49").Value = 'internal static'
$ws.Range("C49").Value = 'public static'
$ws.Range("A50").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxNodeCache.cs'
$ws.Range("B50").Value = 'internal static class SyntaxNodeCache'
$ws.Range("C50").Value = 'public static class SyntaxNodeCache'
$ws.Range("A51").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxNodeCache.cs'
$ws.Range("B51").Value = 'internal static void AddNode'
$ws.Range("C51").Value = 'public static void AddNode'
$ws.Range("A52").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxNodeCache.cs'
$ws.Range("B52").Value = 'internal static GreenNode? TryGetNode(int kind, GreenNode? child1, out int hash)'
$ws.Range("C52").Value = 'public static GreenNode? TryGetNode(int kind, GreenNode? child1, out int hash)'
$ws.Range("A53").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxNodeCache.cs'
$ws.Range("B53").Value = 'internal static GreenNode? TryGetNode(int kind, GreenNode? child1, GreenNode? child2, out int hash)'
$ws.Range("C53").Value = 'public static GreenNode? TryGetNode(int kind, GreenNode? child1, GreenNode? child2, out int hash)'
$ws.Range("A54").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxNodeCache.cs'
$ws.Range("B54").Value = 'internal static GreenNode? TryGetNode(int kind, GreenNode? child1, GreenNode? child2, GreenNode? child3, out int hash)'
$ws.Range("C54").Value = 'public static GreenNode? TryGetNode(int kind, GreenNode? child1, GreenNode? child2, GreenNode? child3, out int hash)'
$ws.Range("A55").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxListPool.cs'
$ws.Range("B55").Value = 'internal'
$ws.Range("C55").Value = 'public'
$ws.Range("A56").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxList`1.cs'
$ws.Range("B56").Value = 'internal partial struct SyntaxList'
$ws.Range("C56").Value = 'public partial struct SyntaxList'
$ws.Range("A57").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxList`1.cs'
$ws.Range("B57").Value = 'internal SyntaxList'
$ws.Range("C57").Value = 'public SyntaxList'
$ws.Range("A58").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxList`1.cs'
$ws.Range("B58").Value = 'internal GreenNode? Node'
$ws.Range("C58").Value = 'public GreenNode? Node'
$ws.Range("A59").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxList`1.Enumerator.cs'
$ws.Range("B59").Value = 'internal partial struct SyntaxList'
$ws.Range("C59").Value = 'public partial struct SyntaxList'
$ws.Range("A60").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxList`1.Enumerator.cs'
$ws.Range("B60").Value = 'internal struct Enumerator'
$ws.Range("C60").Value = 'public struct Enumerator'
$ws.Range("A61").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SeparatedSyntaxList.cs'
$ws.Range("B61").Value = 'internal struct SeparatedSyntaxList'
$ws.Range("C61").Value = 'public struct SeparatedSyntaxList'
$ws.Range("A62").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SeparatedSyntaxList.cs'
$ws.Range("B62").Value = 'internal SeparatedSyntaxList'
$ws.Range("C62").Value = 'public SeparatedSyntaxList'
$ws.Range("A63").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SeparatedSyntaxList.cs'
$ws.Range("B63").Value = 'internal GreenNode? Node'
$ws.Range("C63").Value = 'public GreenNode? Node'
$ws.Range("A64").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxListBuilder.cs'
$ws.Range("B64").Value = 'internal class SyntaxListBuilder'
$ws.Range("C64").Value = 'public class SyntaxListBuilder'
$ws.Range("A65").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SyntaxListBuilder`1.cs'
$ws.Range("B65").Value = 'internal struct SyntaxListBuilder'
$ws.Range("C65").Value = 'public struct SyntaxListBuilder'
$ws.Range("A66").Value = 'MetaDslx.CodeAnalysis.Common\Syntax\InternalSyntax\SeparatedSyntaxListBuilder.cs'
$ws.Range("B66").Value = 'internal struct SeparatedSyntaxListBuilder'
$ws.Range("C66").Value = 'public struct SeparatedSyntaxListBuilder'
$ws.Range("A67").Value = 'MetaDslx.CodeAnalysis.Common\Symbols\TypeCompareKind.cs'
$ws.Range("B67").Value = 'internal enum TypeCompareKind'
$ws.Range("C67").Value = 'public enum TypeCompareKind'
$ws.Range("A68").Value = 'MetaDslx.CodeAnalysis.Common\Symbols\Attributes\ObsoleteAttributeKind.cs'
$ws.Range("B68").Value = 'internal enum ObsoleteAttributeKind'
$ws.Range("C68").Value = 'public enum ObsoleteAttributeKind'
$ws.Range("A69").Value = 'MetaDslx.CodeAnalysis.Common\Symbols\Attributes\ObsoleteAttributeKind.cs'
$ws.Range("B69").Value = 'internal sealed class ObsoleteAttributeData'
$ws.Range("C69").Value = 'public sealed class ObsoleteAttributeData'
$ws.Range("A70").Value = 'MetaDslx.CodeAnalysis.Common\MetadataReader\MetadataTypeName.cs'
$ws.Range("B70").Value = 'internal partial struct MetadataTypeName'
$ws.Range("C70").Value = 'public partial struct MetadataTypeName'
$ws.Range("A71").Value = 'MetaDslx.CodeAnalysis.Common\MetadataReader\MetadataTypeName.Key.cs'
$ws.Range("B71").Value = 'internal partial struct MetadataTypeName'
$ws.Range("C71").Value = 'public partial struct MetadataTypeName'
$ws.Range("A72").Value = 'MetaDslx.CodeAnalysis.Common\MetadataReader\PEModule.cs'
$ws.Range("B72").Value = 'internal sealed class PEModule'
$ws.Range("C72").Value = 'public sealed class PEModule'
$ws.Range("A73").Value = 'MetaDslx.CodeAnalysis.Common\ReferenceManager\CommonReferenceManager.State.cs'
$ws.Range("B73").Value = 'internal abstract class CommonReferenceManager'
$ws.Range("C73").Value = 'public abstract class CommonReferenceManager'
$ws.Range("A74").Value = 'MetaDslx.CodeAnalysis.Common\Compilation\SemanticModelProvider.cs'
$ws.Range("B74").Value = 'internal abstract class SemanticModelProvider'
$ws.Range("C74").Value = 'public abstract class SemanticModelProvider'
$ws.Range("A75").Value = 'MetaDslx.CodeAnalysis.Common\Compilation\Compilation.cs'
$ws.Range("B75").Value = 'internal static void CheckSubmissionOptions'
$ws.Range("C75").Value = 'protected static void CheckSubmissionOptions'
$ws.Range("A76").Value = 'MetaDslx.CodeAnalysis.Common\Compilation\Compilation.cs'
$ws.Range("B76").Value = 'internal static void ValidateScriptCompilationParameters'
$ws.Range("C76").Value = 'protected static void ValidateScriptCompilationParameters'

$ws.Range("C9").Select()
